# The commit inserts one new weekly price record for "Ciboulette" just
# before what is currently row 233, pushing the existing rows 233:368
# down to 234:369 (dimension grows from A1:R368 to A1:R369).
#
# Insert a new row at position 233; Excel shifts rows 233:368 -> 234:369
# automatically (carrying over all of their data/formatting), and the
# new blank row 233 inherits the column D date-number-format style from
# the row that used to occupy that position.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(233).Insert()

# Populate the newly inserted row 233 with the new record's values.
$ws.Cells.Item(233, 1).Value = 3
$ws.Cells.Item(233, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(233, 3).Value = "Coquimbo"
$ws.Cells.Item(233, 4).Value = 44830
$ws.Cells.Item(233, 5).Value = 5
$ws.Cells.Item(233, 6).Value = 100112039
$ws.Cells.Item(233, 7).Value = "Ciboulette"
$ws.Cells.Item(233, 8).Value = "Sin especificar"
$ws.Cells.Item(233, 9).Value = "Primera"
$ws.Cells.Item(233, 10).Value = 120
$ws.Cells.Item(233, 11).Value = 1500
$ws.Cells.Item(233, 12).Value = 1500
$ws.Cells.Item(233, 13).Value = 1500
$ws.Cells.Item(233, 14).Value = "`$/docena de atados"
$ws.Cells.Item(233, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(233, 16).Value = 500
$ws.Cells.Item(233, 17).Value = 3
$ws.Cells.Item(233, 18).Value = "Hortaliza"
